$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04271373187048222
$ws.Range("C2").Value = 0.306821227259698
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("G2").Value = 1155.776987104807
